# Error Calculations and Plots
# This workbook holds a "missing data" table (ID, A, B, C, D, F columns).
# The edit:
#   1) Removes two records entirely (rows whose ID is "RM 232" and "SC 92"),
#      which shifts every following row up.
#   2) After the shift, a number of individual cells in columns D/E/F are
#      updated: some previously-missing cells get a restored numeric value,
#      and a few previously-filled cells become missing (cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the two removed records -------------------------------------
# "RM 232" is row 26 in the original layout.
$ws.Rows(26).Delete()
# After that deletion, "SC 92" (originally row 28) is now row 27.
$ws.Rows(27).Delete()

# --- 2) Apply the per-cell value changes (row numbers below are the final,
#        post-deletion row numbers) -----------------------------------------
$ws.Range("F2").ClearContents()

$ws.Range("F5").Value = 17.66

$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43

$ws.Range("E8").ClearContents()

$ws.Range("F10").ClearContents()

$ws.Range("E12").Value = -5.3

$ws.Range("F13").ClearContents()

$ws.Range("E14").ClearContents()

$ws.Range("E17").Value = -7.3

$ws.Range("E18").Value = -8.5

$ws.Range("E19").ClearContents()

$ws.Range("E20").ClearContents()

$ws.Range("E23").Value = -7

$ws.Range("F24").Value = 16.78

$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()

$ws.Range("F28").ClearContents()

$ws.Range("D29").ClearContents()

$ws.Range("F30").Value = 16.89

$ws.Range("D32").ClearContents()
